$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 2.6 = 9641.3 pesos"), "1000 Bs = 2.63 = 9745.37 pesos"
$newText = $newText -replace [regex]::Escape("9641.3 pesos = 2.59 = 946.61 Bs"), "9745.37 pesos = 2.61 = 954.03 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 380.9
$wsTasas.Range("O10").Value = 3712.01
$wsTasas.Range("O12").Value = 365
